# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# Update OFF sheet (row 2: Week 14 -> Week 15 values)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 343
$wsOff.Range("C2").Value = 218
$wsOff.Range("D2").Value = 82
$wsOff.Range("E2").Value = 31
$wsOff.Range("F2").Value = 10
$wsOff.Range("G2").Value = 5

# Update DEF sheet (row 2: Week 14 -> Week 15 values)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 355
$wsDef.Range("C2").Value = 263
$wsDef.Range("D2").Value = 80
$wsDef.Range("E2").Value = 36
